$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'28.981.52"
$ws.Range("E2").Value = "  -1.98%  "
$ws.Range("D3").Value = "'1.905.00"
$ws.Range("E3").Value = "  -4.32%  "
$ws.Range("E4").Value = "  -0.41%  "
$ws.Range("D5").Value = "'324.06"
$ws.Range("E5").Value = "  -0.97%  "
$ws.Range("D6").Value = "'1.001"
$ws.Range("E6").Value = "  -0.55%  "
$ws.Range("D7").Value = "'0.4596"
$ws.Range("E7").Value = "  -1.85%  "
$ws.Range("D8").Value = "'0.3829"
$ws.Range("E8").Value = "  -3.02%  "
$ws.Range("D9").Value = "'0.07750"
$ws.Range("E9").Value = "  -2.76%  "
$ws.Range("D10").Value = "'0.9816"
$ws.Range("E10").Value = "  -2.21%  "
$ws.Range("D11").Value = "'22.14"
$ws.Range("E11").Value = "  -3.06%  "
$ws.Range("D12").Value = "'1.924.07"
$ws.Range("E12").Value = "  -5.61%  "
$ws.Range("D13").Value = "'6.992"
$ws.Range("E13").Value = "  -3.73%  "
$ws.Range("D14").Value = "'5.696"
$ws.Range("E14").Value = "  -3.24%  "
$ws.Range("D15").Value = "'0.07043"
$ws.Range("E15").Value = "  -1.60%  "
$ws.Range("D16").Value = "'1.003"
$ws.Range("E16").Value = "  -0.46%  "
$ws.Range("D17").Value = "'83.97"
$ws.Range("E17").Value = "  -5.60%  "
$ws.Range("D18").Value = "'0.000009546"
$ws.Range("E18").Value = "  -4.47%  "
$ws.Range("D19").Value = "'16.72"
$ws.Range("E19").Value = "  -3.58%  "
$ws.Range("D20").Value = "'1.001"
$ws.Range("E20").Value = "  -0.64%  "
$ws.Range("D21").Value = "'28.987.63"
$ws.Range("E21").Value = "  -2.25%  "
$ws.Range("D22").Value = "'5.327"
$ws.Range("E22").Value = "  -4.00%  "
$ws.Range("D23").Value = "'10.96"
$ws.Range("E23").Value = "  -2.97%  "
$ws.Range("D24").Value = "'2.125.68"
$ws.Range("E24").Value = "  -6.16%  "
$ws.Range("D25").Value = "'2.071"
$ws.Range("E25").Value = "  -2.71%  "
$ws.Range("D26").Value = "'156.13"
$ws.Range("E26").Value = "  -1.54%  "
$ws.Range("D27").Value = "'19.14"
$ws.Range("E27").Value = "  -2.81%  "
$ws.Range("D28").Value = "'5.625"
$ws.Range("E28").Value = "  -5.68%  "
$ws.Range("D29").Value = "'117.72"
$ws.Range("E29").Value = "  -2.29%  "
$ws.Range("D30").Value = "'1.830"
$ws.Range("E30").Value = "  -7.01%  "
$ws.Range("D31").Value = "'0.09262"
$ws.Range("E31").Value = "  -1.94%  "
$ws.Range("D32").Value = "'0.8601"
$ws.Range("E32").Value = "  -3.87%  "
$ws.Range("D33").Value = "'5.110"
$ws.Range("E33").Value = "  -3.53%  "
$ws.Range("D34").Value = "'1.249"
$ws.Range("E34").Value = "  -7.31%  "
$ws.Range("D35").Value = "'3.014"
$ws.Range("E35").Value = "  -5.69%  "
$ws.Range("D36").Value = "'0.05729"
$ws.Range("E36").Value = "  -1.93%  "
$ws.Range("D37").Value = "'1.151"
$ws.Range("E37").Value = "  -2.25%  "
$ws.Range("D38").Value = "'1.001"
$ws.Range("E38").Value = "  -0.76%  "
$ws.Range("D39").Value = "'0.02048"
$ws.Range("E39").Value = "  -3.99%  "
$ws.Range("D40").Value = "'7.468"
$ws.Range("E40").Value = "  -5.85%  "
$ws.Range("D41").Value = "'0.5529"
$ws.Range("E41").Value = "  -4.11%  "
$ws.Range("D42").Value = "'0.1758"
$ws.Range("E42").Value = "  -3.57%  "
$ws.Range("D43").Value = "'9.284"
$ws.Range("E43").Value = "  -5.73%  "
$ws.Range("B44").Value = "MXToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D44").Value = "'2.720"
$ws.Range("E44").Value = "  +2.94%  "
$ws.Range("B45").Value = "PEPE"
$ws.Range("C45").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D45").Value = "'0.000002800"
$ws.Range("E45").Value = "  -9.42%  "
$ws.Range("D46").Value = "'0.5214"
$ws.Range("E46").Value = "  -3.20%  "
$ws.Range("D47").Value = "'11.28"
$ws.Range("E47").Value = "  -6.65%  "
$ws.Range("D48").Value = "'2.099"
$ws.Range("E48").Value = "  -3.05%  "
$ws.Range("D49").Value = "'0.06827"
$ws.Range("E49").Value = "  -2.10%  "
$ws.Range("D50").Value = "'111.89"
$ws.Range("E50").Value = "  -2.56%  "
$ws.Range("D51").Value = "'1.783"
$ws.Range("E51").Value = "  -4.71%  "
